$wb = $excel.ActiveWorkbook

# The workbook has two sheets ("展览" and "全部类型") that contain identical
# data tables; both need the same "想去人数" (F column) counts updated.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1664
    $ws.Range("F6").Value = 439
    $ws.Range("F8").Value = 68
    $ws.Range("F9").Value = 561
}
